$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date text in column A (A2:A6 all share the "2025-11-27" string).
# Force text format first so Excel doesn't auto-convert the date-like string
# into a date serial number, then restore the default "Normal" style so the
# cell formatting stays identical to before the edit.
$dateRng = $ws.Range("A2:A6")
$dateRng.NumberFormat = "@"
$dateRng.Value = "2025-11-28"
$dateRng.Style = "Normal"

# Update the refreshed metrics for the Bitcoin (row 2) record
$ws.Range("D2").Value = 90801.8
$ws.Range("E2").Value = 33.6
$ws.Range("F2").Value = 7.27
$ws.Range("J2").Value = 40
